# Rename stimuli image identifiers in columns F (firstScore) and G (secondScore)
# for data rows 2-9: a_e1/a_e2 -> stim_a, t_e1/t_e2 -> stim_t.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = "stim_a"
$ws.Range("G2").Value = "stim_t"
$ws.Range("F3").Value = "stim_a"
$ws.Range("G3").Value = "stim_t"
$ws.Range("F4").Value = "stim_a"
$ws.Range("G4").Value = "stim_t"
$ws.Range("F5").Value = "stim_a"
$ws.Range("G5").Value = "stim_t"
$ws.Range("F6").Value = "stim_t"
$ws.Range("G6").Value = "stim_a"
$ws.Range("F7").Value = "stim_t"
$ws.Range("G7").Value = "stim_a"
$ws.Range("F8").Value = "stim_t"
$ws.Range("G8").Value = "stim_a"
$ws.Range("F9").Value = "stim_t"
$ws.Range("G9").Value = "stim_a"

# Match the author's final selection before saving.
$ws.Range("G6:G9").Select()
